$wb = $excel.ActiveWorkbook

# Sheet ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1350
$ws.Range("I28").Value = 1350
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1350
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -865

# Sheet ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7033.8667
$ws.Range("I62").Value = 4663.5
$ws.Range("J62").Value = 9742.857
$ws.Range("K62").Value = 4663.5
$ws.Range("L62").Value = 9742.857
$ws.Range("M62").Value = -4039.5
$ws.Range("N62").Value = -10990.857

# Sheet ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7033.8667
$ws.Range("I65").Value = 4663.5
$ws.Range("J65").Value = 9742.857
$ws.Range("K65").Value = 23317.5
$ws.Range("L65").Value = 48714.285
$ws.Range("M65").Value = -20197.5
$ws.Range("N65").Value = -54954.285

# Sheet ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 64649.25
$ws.Range("I70").Value = 5443.6665
$ws.Range("J70").Value = 84384.44500000001
$ws.Range("K70").Value = 16330.9995
$ws.Range("L70").Value = 253153.335
$ws.Range("M70").Value = -16060.9995
$ws.Range("N70").Value = -253693.335

# Sheet ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 64649.25
$ws.Range("I73").Value = 5443.6665
$ws.Range("J73").Value = 84384.44500000001
$ws.Range("K73").Value = 16330.9995
$ws.Range("L73").Value = 253153.335
$ws.Range("M73").Value = -15394.9995
$ws.Range("N73").Value = -255025.335

# Sheet ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5500
$ws.Range("I116").Value = 5000
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -1558
$ws.Range("N116").Value = -13884

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1707
$ws.Range("I132").Value = 1707
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5121
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2591
$ws.Range("N132").ClearContents()

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3250.2917
$ws.Range("I137").Value = 1498.5
$ws.Range("J137").Value = 4501.5713
$ws.Range("K137").Value = 4495.5
$ws.Range("L137").Value = 13504.7139
$ws.Range("M137").Value = -1945.5

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4160.304
$ws.Range("I138").Value = 1799.1111
$ws.Range("J138").Value = 4734.6484
$ws.Range("K138").Value = 5397.3333
$ws.Range("L138").Value = 14203.9452
$ws.Range("M138").Value = -257.3333000000002
$ws.Range("N138").Value = -24483.9452

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7485.884
$ws.Range("I32").Value = 5689.564
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 5689.564
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -5402.564

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1655.3334
$ws.Range("I45").Value = 1655.3334
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1655.3334
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1278.3334
$ws.Range("N45").ClearContents()

# Sheet ARM row 96
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 6781
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 6781
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 6781
$ws.Range("N96").Value = -12273

# Sheet ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3259.2856
$ws.Range("I97").Value = 1718.3334
$ws.Range("J97").Value = 12505
$ws.Range("K97").Value = 1718.3334
$ws.Range("L97").Value = 12505
$ws.Range("M97").Value = -1222.3334
$ws.Range("N97").Value = -13497

# Sheet ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4360.625
$ws.Range("I110").Value = 4891.4287
$ws.Range("J110").Value = 645
$ws.Range("K110").Value = 4891.4287
$ws.Range("L110").Value = 645
$ws.Range("M110").Value = -2846.4287

# Sheet BSM row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 23666.666
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 23666.666
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 23666.666
$ws.Range("N103").Value = -26010.666

# Sheet BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3002.9092
$ws.Range("I105").Value = 2199.9033
$ws.Range("J105").Value = 15449.5
$ws.Range("K105").Value = 2199.9033
$ws.Range("L105").Value = 15449.5
$ws.Range("M105").Value = -452.9032999999999

# Sheet BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 729.45
$ws.Range("I107").Value = 734.7059
$ws.Range("J107").Value = 699.6667
$ws.Range("K107").Value = 734.7059
$ws.Range("L107").Value = 699.6667
$ws.Range("M107").Value = 1185.2941

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2990.5833
$ws.Range("I134").Value = 2430.5557
$ws.Range("J134").Value = 4670.6665
$ws.Range("K134").Value = 7291.6671
$ws.Range("L134").Value = 14011.9995
$ws.Range("M134").Value = -4756.6671

# Sheet CRP row 28
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 13026
$ws.Range("I28").Value = 10000
$ws.Range("J28").Value = 13530.333
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 13530.333
$ws.Range("M28").Value = -9755
$ws.Range("N28").Value = -14020.333

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4185.074
$ws.Range("I31").Value = 2729.2778
$ws.Range("J31").Value = 7096.6665
$ws.Range("K31").Value = 2729.2778
$ws.Range("L31").Value = 7096.6665
$ws.Range("M31").Value = -2434.2778
$ws.Range("N31").Value = -7686.6665

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4185.074
$ws.Range("I34").Value = 2729.2778
$ws.Range("J34").Value = 7096.6665
$ws.Range("K34").Value = 2729.2778
$ws.Range("L34").Value = 7096.6665
$ws.Range("M34").Value = -2527.2778
$ws.Range("N34").Value = -7500.6665

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5570.125
$ws.Range("I58").Value = 6374.25
$ws.Range("J58").Value = 5302.0835
$ws.Range("K58").Value = 6374.25
$ws.Range("L58").Value = 5302.0835
$ws.Range("M58").Value = -6171.25
$ws.Range("N58").Value = -5708.0835

# Sheet CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 87558.39999999999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 87558.39999999999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 87558.39999999999
$ws.Range("N62").Value = -88806.39999999999

# Sheet CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 87558.39999999999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 87558.39999999999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 437792
$ws.Range("N65").Value = -444032

# Sheet CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8563
$ws.Range("I86").Value = 7232.3335
$ws.Range("J86").Value = 10159.8
$ws.Range("K86").Value = 7232.3335
$ws.Range("L86").Value = 10159.8
$ws.Range("M86").Value = -6109.3335
$ws.Range("N86").Value = -12405.8

# Sheet CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 8563
$ws.Range("I89").Value = 7232.3335
$ws.Range("J89").Value = 10159.8
$ws.Range("K89").Value = 36161.6675
$ws.Range("L89").Value = 50799
$ws.Range("M89").Value = -30545.6675
$ws.Range("N89").Value = -62031

# Sheet CRP row 92
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 29300.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 29300.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 29300.5
$ws.Range("N92").Value = -34292.5

# Sheet CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 49999.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 49999.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 49999.5
$ws.Range("N95").Value = -55491.5

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3284.1304
$ws.Range("I122").Value = 3107.9473
$ws.Range("J122").Value = 4121
$ws.Range("K122").Value = 9323.841899999999
$ws.Range("L122").Value = 12363
$ws.Range("M122").Value = -6873.841899999999
$ws.Range("N122").Value = -17263

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5570.125
$ws.Range("I136").Value = 6374.25
$ws.Range("J136").Value = 5302.0835
$ws.Range("K136").Value = 19122.75
$ws.Range("L136").Value = 15906.2505
$ws.Range("M136").Value = -16572.75
$ws.Range("N136").Value = -21006.2505

# Sheet CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 39.25
$ws.Range("I19").Value = 19.333334
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = 58.000002
$ws.Range("L19").Value = 297
$ws.Range("M19").Value = 115.999998
$ws.Range("N19").Value = -645

# Sheet CUL row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 400273.8
$ws.Range("I29").Value = 666707
$ws.Range("J29").Value = 624
$ws.Range("K29").Value = 2000121
$ws.Range("L29").Value = 1872
$ws.Range("M29").Value = -1999844
$ws.Range("N29").Value = -2426

# Sheet CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 539.4
$ws.Range("I107").Value = 701.5
$ws.Range("J107").Value = 527.8214
$ws.Range("K107").Value = 2104.5
$ws.Range("L107").Value = 1583.4642
$ws.Range("M107").Value = -184.5
$ws.Range("N107").Value = -5423.4642

# Sheet CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1566.6666
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 1720
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 5160
$ws.Range("M121").Value = -1090
$ws.Range("N121").Value = -7780

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1891.1875
$ws.Range("I132").Value = 1269.3572
$ws.Range("J132").Value = 6244
$ws.Range("K132").Value = 3808.0716
$ws.Range("L132").Value = 18732
$ws.Range("M132").Value = -1278.0716
$ws.Range("N132").Value = -23792

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1665.3334
$ws.Range("I61").Value = 1475.8
$ws.Range("J61").Value = 1902.25
$ws.Range("K61").Value = 1475.8
$ws.Range("L61").Value = 1902.25
$ws.Range("M61").Value = -1273.8
$ws.Range("N61").Value = -2306.25

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1665.3334
$ws.Range("I113").Value = 1475.8
$ws.Range("J113").Value = 1902.25
$ws.Range("K113").Value = 1475.8
$ws.Range("L113").Value = 1902.25
$ws.Range("M113").Value = 694.2
$ws.Range("N113").Value = -6242.25

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5003
$ws.Range("I132").Value = 5003
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15009
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12479

# Sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2143.0557
$ws.Range("I122").Value = 2411.7334
$ws.Range("J122").Value = 799.6667
$ws.Range("K122").Value = 7235.2002
$ws.Range("L122").Value = 2399.0001
$ws.Range("M122").Value = -4785.2002

# Sheet WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1957.4
$ws.Range("I126").Value = 1661.4117
$ws.Range("J126").Value = 2586.375
$ws.Range("K126").Value = 4984.2351
$ws.Range("L126").Value = 7759.125
$ws.Range("M126").Value = -2514.2351
$ws.Range("N126").Value = -12699.125

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2214.7144
$ws.Range("I132").Value = 1917.1666
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5751.4998
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3221.4998
$ws.Range("N132").Value = -17060
